$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 20.15985084044064
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 19.48425592650926
$ws.Range("B4").Value = 0.1169995834814548
$ws.Range("C4").Value = 0.3048912486333797
$ws.Range("D4").Value = 3.223369029078222
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 4.178645819794754
$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 3.223369029078222
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.656069925401464
$ws.Range("B6").Value = 0.2881169905109251
$ws.Range("C6").Value = 0.3048912486333797
$ws.Range("D6").Value = 18.71679738969934
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 19.84319158744534
$ws.Range("B7").Value = 0.6545652718822623
$ws.Range("C7").Value = 0.3048912486333797
$ws.Range("D7").Value = 0.7210945179870265
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.213936997104367
$ws.Range("B8").Value = 0.1169995834814548
$ws.Range("C8").Value = 0.3048912486333797
$ws.Range("D8").Value = 3.223369029078222
$ws.Range("E8").Value = 0.5333859586016987
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.178645819794754
$ws.Range("B9").Value = 0.1169995834814548
$ws.Range("C9").Value = 0.3048912486333797
$ws.Range("D9").Value = 18.71679738969934
$ws.Range("E9").Value = 0.5333859586016987
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 19.67207418041587
$ws.Range("B10").Value = 3.272327238179451
$ws.Range("C10").Value = 1.626987699542094
$ws.Range("D10").Value = 0.7210945179870265
$ws.Range("E10").Value = 0.5333859586016987
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.15379541431027
$ws.Range("B11").Value = 1.445647641019636
$ws.Range("C11").Value = 1.626987699542094
$ws.Range("D11").Value = 3.223369029078222
$ws.Range("E11").Value = 0.5333859586016987
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6.82939032824165
$ws.Range("B12").Value = 3.272327238179451
$ws.Range("C12").Value = 1.626987699542094
$ws.Range("D12").Value = 3.223369029078222
$ws.Range("E12").Value = 0.5333859586016987
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 8.656069925401464
